# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 535
    $ws.Range("F4").Value = 269
    $ws.Range("F7").Value = 751
}
